$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values per repulled/recalculated data
$ws.Range("F2").Value = -3
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = -5
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = -2
$ws.Range("F14").Value = -3
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = -5
$ws.Range("F18").Value = 5
$ws.Range("F19").Value = -5
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = -3
$ws.Range("F22").Value = 1
$ws.Range("F23").Value = 7
$ws.Range("F26").Value = -2
$ws.Range("F28").Value = -2
$ws.Range("F29").Value = 1
$ws.Range("F30").Value = -1
$ws.Range("F31").Value = 3
$ws.Range("F33").Value = -8
$ws.Range("F34").Value = 3
$ws.Range("F35").Value = -3
